## Groove libraries.xlsx update
## - "Apache Batik" description cell (libraries!B3) had its trailing
##   ellipsis character replaced by three literal dots ("..."), and the
##   "save to libraries.csv" sheet picks this up automatically through its
##   =libraries!B3 formula.
## - The workbook is left with the "libraries" sheet active/selected (cell
##   B3 selected there), while the "save to libraries.csv" sheet keeps B2
##   as its last selection.
## - The "save to libraries.csv" sheet gains explicit column widths for
##   columns A and B.

$wb = $excel.ActiveWorkbook

$libraries = $wb.Worksheets.Item("libraries")
$saveSheet = $wb.Worksheets.Item("save to libraries.csv")

# Fix the typographic ellipsis to three plain dots.
$libraries.Range("B3").Value = "exporting images to SVG, EPS, ..."

# Give the "save to libraries.csv" sheet explicit column widths (A=30, B=39
# once Excel's +5/6 character-width padding is taken into account).
$saveSheet.Columns.Item(1).ColumnWidth = 29.166666666666668
$saveSheet.Columns.Item(2).ColumnWidth = 38.166666666666664

# Leave "save to libraries.csv" with B2 selected (no longer the active tab) ...
$saveSheet.Range("B2").Select()

# ... then switch to / select the "libraries" sheet, with B3 selected there.
$libraries.Activate()
$libraries.Range("B3").Select()
